# Update the "Förändrad" (Changed) date column (C) on the data rows:
# every occurrence of the old date serial 45203 (2023-10-04) becomes
# 45205 (2023-10-06), for rows 2 through 123.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 123
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45205
    }
}
